$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated td_sim_1 (column C) and record_atd (column D) values for rows 2-16,
# plus the recalculated average in C17.

$ws.Range("C2").Value = 222
$ws.Range("D2").Value = 206.5

$ws.Range("C3").Value = 124
$ws.Range("D3").Value = 99.5

$ws.Range("C4").Value = 111
$ws.Range("D4").Value = 85.5

$ws.Range("C5").Value = 63
$ws.Range("D5").Value = 33.5

$ws.Range("C6").Value = 223
$ws.Range("D6").Value = 207.5

$ws.Range("C7").Value = 80
$ws.Range("D7").Value = 56.5

$ws.Range("C8").Value = 76
$ws.Range("D8").Value = 48.5

$ws.Range("C9").Value = 1446
$ws.Range("D9").Value = 1495.5

$ws.Range("C10").Value = 75
$ws.Range("D10").Value = 43.5

$ws.Range("C11").Value = 118
$ws.Range("D11").Value = 89.5

$ws.Range("C12").Value = 97
$ws.Range("D12").Value = 71

$ws.Range("C13").Value = 94
$ws.Range("D13").Value = 94

$ws.Range("C14").Value = 710
$ws.Range("D14").Value = 813.5

$ws.Range("C15").Value = 64
$ws.Range("D15").Value = 34.5

$ws.Range("C16").Value = 65
$ws.Range("D16").Value = 35.5

$ws.Range("C17").Value = 237.8666666666667
